# Reorders the "Artfynd" sheet's species-observation rows.
#
# The underlying edit is a pure row-content permutation: whole rows
# (17..28 and 31,33,34,35) had their entire contents (all of columns
# A:AY) swapped around between each other, while the row numbers
# themselves (and everything outside those rows, e.g. headers / row 29,
# 30, 32) stayed put.
#
# Because it is a permutation made up of disjoint cycles, every row in a
# cycle must be read into memory BEFORE any of them are written back,
# otherwise an earlier write would clobber data that a later read still
# needs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last used column is AY (=51).
$maxCol = 51

function Read-RowData {
    param($ws, $row, $maxCol)
    $data = @{}
    for ($c = 1; $c -le $maxCol; $c++) {
        $data[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $data
}

function Write-RowData {
    param($ws, $row, $data, $maxCol)
    for ($c = 1; $c -le $maxCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $data[$c]
    }
}

# Disjoint row-permutation cycles: for each row in a cycle, its new
# content comes from the NEXT row listed (wrapping around).
#   new(17) = old(18), new(18) = old(17)
#   new(19) = old(26), new(26) = old(20), new(20) = old(22), new(22) = old(19)
#   new(21) = old(23), new(23) = old(24), new(24) = old(25), new(25) = old(21)
#   new(27) = old(28), new(28) = old(27)
#   new(31) = old(34), new(34) = old(31)
#   new(33) = old(35), new(35) = old(33)
$cycles = @(
    @(17, 18),
    @(19, 26, 20, 22),
    @(21, 23, 24, 25),
    @(27, 28),
    @(31, 34),
    @(33, 35)
)

foreach ($cycle in $cycles) {
    # 1) Snapshot every row in this cycle first.
    $snapshots = @{}
    foreach ($r in $cycle) {
        $snapshots[$r] = Read-RowData $ws $r $maxCol
    }

    # 2) Write each row's new content: row at position i gets the data
    #    that was read from the row at position i+1 (wrapping around).
    $n = $cycle.Count
    for ($i = 0; $i -lt $n; $i++) {
        $targetRow = $cycle[$i]
        $sourceRow = $cycle[($i + 1) % $n]
        Write-RowData $ws $targetRow $snapshots[$sourceRow] $maxCol
    }
}

Write-Host "Row permutation applied."
